$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "skfaycx@erxnbzm.com"
$ws.Range("B1").Value = "H2NyiyK"
$ws.Range("C1").Value = "qsgnucc"
$ws.Range("D1").Value = "Foodie."
